# Apply "Included additional review comments on Overview" edits to slide 2.
#
# The three edited text boxes all live inside a deeply-nested top-level
# group ("Group 1"); PowerPoint's GroupItems collection flattens that
# nesting so every leaf shape can be reached directly from it.
#
# NOTE: this runtime's Shape.Left/.Top/.Width/.Height setters write the
# point value straight into the shape's local <a:off>/<a:ext> (no group
# chOff/chExt re-scaling) after round-tripping the value through a
# float32, so the literals below were chosen so that
# floor(float32(pt) * 12700) reproduces the exact target EMU from the
# OOXML diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$group1 = $s.Shapes.Item(1)
$items = $group1.GroupItems

# --- TextBox 129 (id=130): "Monitors VM " -> "Monitors  ", widen box ---
$monitorsBox = $items.Item("TextBox 129")
$monitorsBox.Left = 190.56536865234375
$monitorsBox.Width = 95.6712646484375
$monitorsBox.TextFrame.TextRange.Text = "Monitors  "

# --- TextBox 42 (id=43): "for Service A " -> "Service A and B ", widen box ---
$serviceBox = $items.Item("TextBox 42")
$serviceBox.Left = 192.1956024169922
$serviceBox.Width = 99.7653579711914
$serviceBox.TextFrame.TextRange.Text = "Service A and B "

# --- TextBox 44 (id=45): reposition/widen only, text unchanged ---
$monitoringSelfBox = $items.Item("TextBox 44")
$monitoringSelfBox.Left = 439.37890625
$monitoringSelfBox.Width = 122.16055297851562
